$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 26 de Junio de 2020 a las 19:42"

# Apply updated COVID-19 country statistics (values + reordered rows)
$ws.Range("B4").Value = 2527919
$ws.Range("C4").Value = 23331
$ws.Range("D4").Value = 1054917
$ws.Range("E4").Value = 1345902
$ws.Range("G4").Value = 320
$ws.Range("H4").Value = 127100
$ws.Range("B5").Value = 1244419
$ws.Range("C5").Value = 11272
$ws.Range("E5").Value = 515386
$ws.Range("G5").Value = 250
$ws.Range("H5").Value = 55304
$ws.Range("B7").Value = 509170
$ws.Range("C7").Value = 18000
$ws.Range("D7").Value = 295917
$ws.Range("E7").Value = 197564
$ws.Range("G7").Value = 381
$ws.Range("H7").Value = 15689
$ws.Range("B8").Value = 309360
$ws.Range("C8").Value = 1380
$ws.Range("G8").Value = 184
$ws.Range("H8").Value = 43414
$ws.Range("A16").Value = "Turquia"
$ws.Range("B16").Value = 194511
$ws.Range("C16").Value = 1396
$ws.Range("D16").Value = 167198
$ws.Range("E16").Value = 22248
$ws.Range("G16").Value = 19
$ws.Range("H16").Value = 5065
$ws.Range("A17").Value = "Alemania"
$ws.Range("B17").Value = 194013
$ws.Range("C17").Value = 228
$ws.Range("D17").Value = 177100
$ws.Range("E17").Value = 7901
$ws.Range("H17").Value = 9012
$ws.Range("B22").Value = 102733
$ws.Range("C22").Value = 111
$ws.Range("D22").Value = 65651
$ws.Range("E22").Value = 28575
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 8507
$ws.Range("B49").Value = 25414
$ws.Range("C49").Value = 9
$ws.Range("E49").Value = 320
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = 1730
$ws.Range("B108").Value = 2014
$ws.Range("C108").Value = 4
$ws.Range("E108").Value = 384
$ws.Range("A119").Value = "Estado de Palestina"
$ws.Range("B119").Value = 1557
$ws.Range("C119").Value = 175
$ws.Range("D119").Value = 446
$ws.Range("E119").Value = 1108
$ws.Range("H119").Value = 3
$ws.Range("A120").Value = "Guinea-Bisau"
$ws.Range("B120").Value = 1556
$ws.Range("C120").Value = 0
$ws.Range("D120").Value = 191
$ws.Range("E120").Value = 1346
$ws.Range("H120").Value = 19
$ws.Range("A121").Value = "Nueva Zelanda"
$ws.Range("B121").Value = 1520
$ws.Range("C121").Value = 1
$ws.Range("D121").Value = 1484
$ws.Range("E121").Value = 14
$ws.Range("H121").Value = 22
$ws.Range("B163").Value = 254
$ws.Range("C163").Value = 12
$ws.Range("E163").Value = 151
$ws.Range("A200").Value = "Laos"
$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("A202").Value = "Dominica"
$ws.Range("A203").Value = "Fiyi"
$ws.Range("A208").Value = "Islas Malvinas"
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A212").Value = "Montserrat"
$ws.Range("D212").Value = 10
$ws.Range("H212").Value = 1
$ws.Range("A213").Value = "Seychelles"
$ws.Range("D213").Value = 11
$ws.Range("H213").Value = 0

Write-Host "Applied country/provincia updates to sheet:" $ws.Name
